$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.654.67"
$ws.Cells.Item(2, 5).Value = "  -0.73%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.504.33"
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "574.38"
$ws.Cells.Item(5, 5).Value = "  -0.67%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "166.52"
$ws.Cells.Item(6, 5).Value = "  -0.95%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 5).Value = "  -1.18%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.502.98"
$ws.Cells.Item(9, 5).Value = "  -0.67%  "
$ws.Cells.Item(11, 5).Value = "  -0.01%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.359"
$ws.Cells.Item(12, 5).Value = "  +2.89%  "
$ws.Cells.Item(13, 5).Value = "  +0.58%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.959.72"
$ws.Cells.Item(14, 5).Value = "  -0.74%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "69.522.52"
$ws.Cells.Item(15, 5).Value = "  -0.73%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.0000177"
$ws.Cells.Item(16, 5).Value = "  +0.52%  "
$ws.Cells.Item(17, 5).Value = "  -1.81%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.512.64"
$ws.Cells.Item(18, 5).Value = "  -0.45%  "
$ws.Cells.Item(19, 5).Value = "  -1.28%  "
$ws.Cells.Item(20, 5).Value = "  -3.75%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "349.29"
$ws.Cells.Item(21, 5).Value = "  -0.63%  "
$ws.Cells.Item(22, 5).Value = "  -0.96%  "
$ws.Cells.Item(23, 5).Value = "  -1.08%  "
$ws.Cells.Item(24, 5).Value = "  -0.10%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "70.84"
$ws.Cells.Item(25, 5).Value = "  +2.16%  "
$ws.Cells.Item(26, 5).Value = "  -2.01%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.75"
$ws.Cells.Item(27, 5).Value = "  -3.26%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.630.95"
$ws.Cells.Item(28, 5).Value = "  -0.66%  "
$ws.Cells.Item(29, 5).Value = "  -0.07%  "
$ws.Cells.Item(30, 5).Value = "  -2.06%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.86"
$ws.Cells.Item(31, 5).Value = "  -0.72%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "458.17"
$ws.Cells.Item(32, 5).Value = "  -1.80%  "
$ws.Cells.Item(33, 5).Value = "  -5.96%  "
$ws.Cells.Item(34, 5).Value = "  -1.62%  "
$ws.Cells.Item(35, 5).Value = "  -0.03%  "
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.116"
$ws.Cells.Item(36, 5).Value = "  -3.11%  "
$ws.Cells.Item(37, 2).Value = "Monero"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "157.00"
$ws.Cells.Item(37, 5).Value = "  +2.47%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "19.06"
$ws.Cells.Item(38, 5).Value = "  +0.18%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "18.35"
$ws.Cells.Item(39, 5).Value = "  -0.75%  "
$ws.Cells.Item(40, 5).Value = "  +0.01%  "
$ws.Cells.Item(41, 5).Value = "  -1.39%  "
$ws.Cells.Item(42, 5).Value = "  -2.04%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.60"
$ws.Cells.Item(43, 5).Value = "  -0.17%  "
$ws.Cells.Item(44, 5).Value = "  -0.07%  "
$ws.Cells.Item(45, 5).Value = "  -5.18%  "
$ws.Cells.Item(46, 5).Value = "  -8.03%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "141.36"
$ws.Cells.Item(47, 5).Value = "  -1.37%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.48"
$ws.Cells.Item(48, 5).Value = "  -0.64%  "
$ws.Cells.Item(50, 5).Value = "  -0.64%  "
$ws.Cells.Item(51, 5).Value = "  -0.49%  "
